$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 26 ("RM 232") - causes rows below to shift up by one
$ws.Rows.Item(26).Delete()

# Delete entire row that is now "SC 92" (was row 28, now row 27 after first delete)
$ws.Rows.Item(27).Delete()

# Now apply individual cell value edits (rows 2-25 unaffected by the shift)
$ws.Range("E2").Value = -7.2
$ws.Range("F3").Value = $null
$ws.Range("F4").Value = 17.97
$ws.Range("E6").Value = $null
$ws.Range("F8").Value = $null
$ws.Range("F9").Value = $null
$ws.Range("E12").Value = -5.3
$ws.Range("E14").Value = $null
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("F19").Value = $null
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("F22").Value = $null
$ws.Range("E23").Value = $null
$ws.Range("F23").Value = 16.48
$ws.Range("E24").Value = $null
$ws.Range("F25").Value = 16.6

# Rows 26-33 (post-deletion) need value corrections relative to the shifted-up data
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = $null
$ws.Range("F27").Value = $null
$ws.Range("B28").Value = $null
$ws.Range("B29").Value = -19.5
$ws.Range("B30").Value = -19.7
$ws.Range("B31").Value = $null
$ws.Range("E31").Value = -8.1
$ws.Range("B32").Value = $null
$ws.Range("E33").Value = -10.7
